# Refactor: extract scenario code.
#
# - Mark rows 85-87 (scenarios #84-#86) as "done" (column B) to match the
#   pattern already used by all earlier rows.
# - Append a new scenario row (row 88, #87) for an XRP backtest that was
#   extracted/split out, re-using the coefficient values that used to sit
#   in row 87.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "done" flag for the existing rows that were missing it.
$ws.Range("B85").Value = "Yes"
$ws.Range("B86").Value = "Yes"
$ws.Range("B87").Value = "Yes"

# Add the new scenario row 88.
$ws.Range("A88").Value = 87
$ws.Range("C88").Value = 200000
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = "XRP"
$ws.Range("F88").Value = "2018-01-06"
$ws.Range("G88").Value = "2018-04-01"
$ws.Range("H88").Value = 0.077137627095034397
$ws.Range("I88").Value = 0.147341631456371
$ws.Range("J88").Value = 0.0033881238556444401
$ws.Range("K88").Value = 0.58523196662204202

# Reflect the user's final selection in the sheet view.
$ws.Range("B84:B87").Select()
